$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename B1 label, add new C1 header describing the new computed column
$ws.Range("B1").Value = "ReDrawing 2500 Rects (w/text) 100 times (sec)"
$ws.Range("C1").Value = "Per frame elapsed time"

# Widen column B to fit the longer header text
$ws.Columns.Item(2).ColumnWidth = 39.43

# Row 2: wxDC (w/ GraphicsContext) -> numeric seconds (was text) + formula dividing by 100
$ws.Range("B2").Value = 17.5099999905
$ws.Range("C2").Formula = "=B2/100"

# Rows 3:4 -> numeric seconds (was text) + shared formula dividing by 100
$ws.Range("B3").Value = 194.18
$ws.Range("B4").Value = 8.81100010872
$ws.Range("C3:C4").Formula = "=B3/100"

# Update the selection to span A1:C4 with active cell C4
$ws.Range("A1:C4").Select()
$ws.Range("C4").Activate()
